# Auto-generated edit script: updates cryptocurrency price (D) and
# 1h volume/change (E) columns to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.786.68'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '3.114.95'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'588.70"
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").Value = "'135.59"
$ws.Range("E6").Value = '  -4.83%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.110.57'
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = '  -4.24%  '
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = '  -3.27%  '
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = '  -5.72%  '
$ws.Range("D14").Value = "'33.71"
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("D15").Value = '3.626.52'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '62.808.31'
$ws.Range("E17").Value = '  -1.90%  '
$ws.Range("D18").Value = '3.111.84'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = "'6.60"
$ws.Range("E19").Value = '  -3.77%  '
$ws.Range("D20").Value = "'466.92"
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").Value = "'13.98"
$ws.Range("E21").Value = '  -3.78%  '
$ws.Range("D22").Value = "'0.690"
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").Value = "'7.59"
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").Value = "'85.11"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = "'12.81"
$ws.Range("E25").Value = '  -4.03%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = '  -6.70%  '
$ws.Range("D29").Value = "'6.79"
$ws.Range("E29").Value = '  -5.08%  '
$ws.Range("D30").Value = "'2.06"
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = "'26.40"
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("E33").Value = '  -6.08%  '
$ws.Range("E34").Value = '  -5.11%  '
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").Value = "'5.69"
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("D38").Value = "0.0$([char]0x2083)0672"
$ws.Range("E38").Value = '  -12.50%  '
$ws.Range("D39").Value = "'0.0383"
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("D40").Value = "'415.39"
$ws.Range("E40").Value = '  -6.40%  '
$ws.Range("D41").Value = "'8.14"
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D42").Value = '2.882.10'
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = '  -12.63%  '
$ws.Range("D44").Value = "'0.111"
$ws.Range("E44").Value = '  -6.27%  '
$ws.Range("D45").Value = "'0.256"
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = "'2.07"
$ws.Range("E47").Value = '  -6.37%  '
$ws.Range("D48").Value = "'25.15"
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").Value = "'2.22"
$ws.Range("E50").Value = '  -8.27%  '
$ws.Range("E51").Value = '  +0.02%  '
